$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Save" header in H1, matching the style of the existing header row (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Value for the new Save column, row 2
$ws.Range("H2").Value = 1
